$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "57.897.84"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "  -0.66%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "2.443.70"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "  -1.36%  "
$ws.Range("E4").Value = "  -0.01%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "511.61"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -1.69%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "129.83"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  -1.62%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.997"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  -0.18%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.550"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  -1.64%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "2.458.71"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  -0.95%  "
$ws.Range("E10").Value = "  -0.01%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.0958"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  -3.85%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "5.19"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  -3.50%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "0.330"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  -4.12%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "2.876.38"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  -1.44%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "57.812.41"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  -0.70%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "21.99"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  -0.61%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "0.0000133"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  -2.55%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "2.448.31"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  -1.36%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "10.50"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  -3.21%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "317.78"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  -1.05%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "4.12"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  -1.76%  "
$ws.Range("E22").Value = "  +0.18%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "5.87"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  +1.56%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "63.13"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  -1.73%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "0.404"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  -1.52%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "0.992"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  -0.74%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "0.159"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  -1.48%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "7.26"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  -2.17%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "168.12"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  +0.90%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "0.0₃0730"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  -3.17%  "
$ws.Range("B31").Value = "PancakeSwap"
$ws.Range("C31").Value = "https://coinranking.com/coin/ncYFcP709+pancakeswap-cake"
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "1.67"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  -2.51%  "
$ws.Range("B32").Value = "Aptos"
$ws.Range("C32").Value = "https://coinranking.com/coin/HGYj5JCv5+aptos-apt"
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "6.20"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  -2.54%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "1.16"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  -3.13%  "
$ws.Range("E34").Value = "  -0.08%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "0.995"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  -0.38%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "17.77"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  -2.25%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "1.27"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  -4.36%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "3.92"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  -2.34%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "36.53"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  -0.21%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "1.46"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  -1.88%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.760"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  -4.65%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "272.08"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  -2.02%  "
$ws.Range("B43").Value = "RenderToken"
$ws.Range("C43").Value = "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "5.01"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  -2.00%  "
$ws.Range("B44").Value = "Filecoin"
$ws.Range("C44").Value = "https://coinranking.com/coin/ymQub4fuB+filecoin-fil"
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "3.38"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  -3.03%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.585"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  -1.99%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.0913"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  +0.31%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "119.68"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  -5.66%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "0.0488"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  -1.10%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "17.23"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  -3.74%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.0210"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  -2.52%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "16.68"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  -3.26%  "
